# Weekly data refresh: insert two new price records (dated 44783) at the
# top of the data block (row 1045) on the single sheet, shifting the
# existing rows (old 1045-1118) down by two (to new 1047-1120). No data is
# lost - the two oldest-in-window rows simply end up at the bottom of the
# sheet (new rows 1119-1120) because nothing is deleted, only inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 1045, shifting everything at/after it
# down by two rows (xlShiftDown = -4121).
$ws.Rows("1045:1046").Insert(-4121)

# Populate the newly inserted row 1045.
$ws.Range("A1045").Value = 10
$ws.Range("B1045").Value = "Vega Modelo de Temuco"
$ws.Range("C1045").Value = "La Araucanía"
$ws.Range("D1045").Value = 44783
$ws.Range("E1045").Value = 9
$ws.Range("F1045").Value = 100112004
$ws.Range("G1045").Value = "Cebolla"
$ws.Range("H1045").Value = "Morada(o)"
$ws.Range("I1045").Value = "1a (guarda)"
$ws.Range("J1045").Value = 80
$ws.Range("K1045").Value = 18000
$ws.Range("L1045").Value = 18000
$ws.Range("M1045").Value = 18000
$ws.Range("N1045").Value = "$/malla 18 kilos"
$ws.Range("O1045").Value = "Región de O'Higgins"
$ws.Range("P1045").Value = 1000
$ws.Range("Q1045").Value = 18
$ws.Range("R1045").Value = "Hortaliza"

# Populate the newly inserted row 1046.
$ws.Range("A1046").Value = 10
$ws.Range("B1046").Value = "Vega Modelo de Temuco"
$ws.Range("C1046").Value = "La Araucanía"
$ws.Range("D1046").Value = 44783
$ws.Range("E1046").Value = 9
$ws.Range("F1046").Value = 100112004
$ws.Range("G1046").Value = "Cebolla"
$ws.Range("H1046").Value = "Sin especificar"
$ws.Range("I1046").Value = "1a (guarda)"
$ws.Range("J1046").Value = 300
$ws.Range("K1046").Value = 7000
$ws.Range("L1046").Value = 7000
$ws.Range("M1046").Value = 7000
$ws.Range("N1046").Value = "$/malla 15 kilos"
$ws.Range("O1046").Value = "Región de O'Higgins"
$ws.Range("P1046").Value = 467
$ws.Range("Q1046").Value = 15
$ws.Range("R1046").Value = "Hortaliza"
